$d = $word.ActiveDocument

$replacements = @(
    @{old="379÷8="; new="288÷3="},
    @{old="794÷7="; new="724÷8="},
    @{old="215÷6="; new="611÷3="},
    @{old="203÷2="; new="956÷6="},
    @{old="176÷6="; new="527÷6="},
    @{old="169÷3="; new="256÷2="},
    @{old="878÷3="; new="804÷9="},
    @{old="326÷4="; new="518÷4="},
    @{old="318÷9="; new="782÷5="},
    @{old="485÷6="; new="377÷7="},
    @{old="193÷5="; new="625÷8="},
    @{old="586÷7="; new="888÷7="},
    @{old="285÷8="; new="298÷4="},
    @{old="451÷2="; new="909÷2="},
    @{old="663÷9="; new="640÷4="},
    @{old="628÷7="; new="611÷5="},
    @{old="321÷5="; new="279÷3="},
    @{old="953÷5="; new="759÷2="},
    @{old="792÷5="; new="914÷9="},
    @{old="577÷4="; new="393÷4="},
    @{old="139÷9="; new="109÷7="},
    @{old="457÷5="; new="219÷9="},
    @{old="704÷6="; new="546÷3="},
    @{old="758÷7="; new="678÷2="},
    @{old="515÷9="; new="250÷6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
